$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$date = "2025-10-07 12:36:51"

# ---- Update the timestamp in column A for every pre-existing data row (2-8) ----
$ws.Range("A2").Value = $date
$ws.Range("A3").Value = $date
$ws.Range("A4").Value = $date
$ws.Range("A5").Value = $date
$ws.Range("A6").Value = $date
$ws.Range("A7").Value = $date
$ws.Range("A8").Value = $date

# ---- Row 3: brand-new job inserted right after row 2 ----
$ws.Range("B3").Value = "自社開発のロジシステムをサポート及び開発できる方募集【PHP, Python, VBA etc】"
$ws.Range("F3").Value = "https://www.lancers.jp/work/detail/5408668"
$ws.Range("G3").Value = 305
$ws.Range("H3").Value = "🔥Python ◆開発 ○PHP"

# ---- Row 4: old row 3's job, shifted down one row ----
$ws.Range("B4").Value = "【医療機関向け】既存システム改修・機能拡充エンジニア(Node.js/TypeScript)"
$ws.Range("D4").Value = "300,000 円 ~ 500,000 円 / 固定"
$ws.Range("F4").Value = "https://www.lancers.jp/work/detail/5408214"
$ws.Range("G4").Value = 185
$ws.Range("H4").Value = "🔥TypeScript ◆Node.js"

# ---- Row 5: old row 4's job, shifted down one row ----
$ws.Range("B5").Value = "【急募】携帯アプリ開発のプロフェッショナルを探しています!"
$ws.Range("D5").Value = "500,000 円 ~ 1,000,000 円 / 固定"
$ws.Range("F5").Value = "https://www.lancers.jp/work/detail/5408519"
$ws.Range("G5").Value = 100
$ws.Range("H5").Value = "◆開発 ◇アプリ"

# ---- Row 6: brand-new job ----
$ws.Range("B6").Value = "2026年度新入社員研修Javaサブ講師 (4~6月)"
$ws.Range("D6").Value = "1,000,000 円 ~ 3,000,000 円 / 固定"
$ws.Range("F6").Value = "https://www.lancers.jp/work/detail/5408522"
$ws.Range("G6").Value = 85
$ws.Range("H6").Value = "★Java"

# ---- Row 7: brand-new job ----
$ws.Range("B7").Value = "2026年度新入社員研修Javaサブ講師"
$ws.Range("D7").Value = "500,000 円 ~ 1,000,000 円 / 固定"
$ws.Range("F7").Value = "https://www.lancers.jp/work/detail/5408524"
$ws.Range("G7").Value = 85
$ws.Range("H7").Value = "★Java"

# ---- Row 8: old row 5's job, shifted down three rows ----
$ws.Range("B8").Value = "2026年度新入社員研修Javaメイン講師"
$ws.Range("D8").Value = "1,000,000 円 ~ 3,000,000 円 / 固定"
$ws.Range("F8").Value = "https://www.lancers.jp/work/detail/5408503"
$ws.Range("G8").Value = 85
$ws.Range("H8").Value = "★Java"

# ---- Row 9: brand-new row ----
$ws.Range("A9").Value = $date
$ws.Range("B9").Value = "【急募】LLMによるMCP(Model Context Protocol)でのExcel操作機能開発"
$ws.Range("C9").Value = "システム開発"
$ws.Range("D9").Value = "500,000 円 ~ 1,000,000 円 / 固定"
$ws.Range("E9").Value = "期限情報なし"
$ws.Range("F9").Value = "https://www.lancers.jp/work/detail/5400689"
$ws.Range("G9").Value = 75
$ws.Range("H9").Value = "◆開発"

# ---- Row 10: brand-new row ----
$ws.Range("A10").Value = $date
$ws.Range("B10").Value = "仮想通貨トレードの運用とコンサル【1名】のみ募集"
$ws.Range("C10").Value = "システム開発"
$ws.Range("D10").Value = "300,000 円 ~ 500,000 円 / 固定"
$ws.Range("E10").Value = "期限情報なし"
$ws.Range("F10").Value = "https://www.lancers.jp/work/detail/5393695"
$ws.Range("G10").Value = 55
$ws.Range("H10").Value = "◆コンサル"

# ---- Row 11: brand-new row ----
$ws.Range("A11").Value = $date
$ws.Range("B11").Value = "【急募】社内Webアプリケーションのセキュリティ診断依頼"
$ws.Range("C11").Value = "システム開発"
$ws.Range("D11").Value = "500,000 円 ~ 1,000,000 円 / 固定"
$ws.Range("E11").Value = "期限情報なし"
$ws.Range("F11").Value = "https://www.lancers.jp/work/detail/5408563"
$ws.Range("G11").Value = 45
$ws.Range("H11").Value = "◇アプリ"

# ---- Row 12: old row 6's job (Notion), shifted down six rows ----
$ws.Range("A12").Value = $date
$ws.Range("B12").Value = "Notion実装・運用管理パートナー募集(長期・リモート)"
$ws.Range("C12").Value = "システム開発"
$ws.Range("D12").Value = "50,000 円 ~ 100,000 円 / 固定"
$ws.Range("E12").Value = "期限情報なし"
$ws.Range("F12").Value = "https://www.lancers.jp/work/detail/5408286"
$ws.Range("G12").Value = 38
$ws.Range("H12").Value = "◇管理"

# ---- Row 13: old row 7's job (iOS/Android subscription), shifted down six rows ----
$ws.Range("A13").Value = $date
$ws.Range("B13").Value = "初回 iOSとAndroidのアプリ 課金(サブスク)"
$ws.Range("C13").Value = "システム開発"
$ws.Range("D13").Value = "10,000 円 ~ 20,000 円 / 固定"
$ws.Range("E13").Value = "期限情報なし"
$ws.Range("F13").Value = "https://www.lancers.jp/work/detail/5398382"
$ws.Range("G13").Value = 30
$ws.Range("H13").Value = "◇アプリ"

# ---- Row 14: old row 8's job (Ctrl-click folder app), shifted down six rows ----
$ws.Range("A14").Value = $date
$ws.Range("B14").Value = "Ctrlを押しながらフォルダの上をクリックすると別窓で上の階層のフォルダが開くアプリの作成"
$ws.Range("C14").Value = "システム開発"
$ws.Range("D14").Value = "10,000 円 ~ 20,000 円 / 固定"
$ws.Range("E14").Value = "期限情報なし"
$ws.Range("F14").Value = "https://www.lancers.jp/work/detail/5408148"
$ws.Range("G14").Value = 30
$ws.Range("H14").Value = "◇アプリ"

# ---- Row 15: brand-new row (no H value) ----
$ws.Range("A15").Value = $date
$ws.Range("B15").Value = "【急募】Oracleを活用したQ&Aシステム構築の依頼"
$ws.Range("C15").Value = "システム開発"
$ws.Range("D15").Value = "300,000 円 ~ 500,000 円 / 固定"
$ws.Range("E15").Value = "期限情報なし"
$ws.Range("F15").Value = "https://www.lancers.jp/work/detail/5408554"
$ws.Range("G15").Value = 40

# ---- Row 16: brand-new row (no H value) ----
$ws.Range("A16").Value = $date
$ws.Range("B16").Value = "急募 限定公開 限定公開の仕事"
$ws.Range("C16").Value = "システム開発"
$ws.Range("D16").Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Range("E16").Value = "期限情報なし"
$ws.Range("F16").Value = "https://www.lancers.jp/work/detail/5408546"
$ws.Range("G16").Value = 13

# ---- Rebuild every hyperlink in column F (rows 2-16) ----
# The 7 pre-existing hyperlinks (F2:F8) were bound to URLs that now belong to
# different rows after the reshuffle above, and this engine's Hyperlink
# object offers no working in-place retarget (Address re-assignment and
# Item(n).Delete() are both no-ops here), while Range.Hyperlinks.Delete()
# clears hyperlinks sheet-wide rather than per-cell. So: clear everything
# once, then re-Add all 15 links in row order so the relationship ids come
# out sequential again (rId1->F2, rId2->F3, ... rId15->F16), matching the
# final ref/r:id layout.
$ws.Range("F2").Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("F2"), "https://www.lancers.jp/work/detail/5408156") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F3"), "https://www.lancers.jp/work/detail/5408668") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F4"), "https://www.lancers.jp/work/detail/5408214") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F5"), "https://www.lancers.jp/work/detail/5408519") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F6"), "https://www.lancers.jp/work/detail/5408522") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F7"), "https://www.lancers.jp/work/detail/5408524") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F8"), "https://www.lancers.jp/work/detail/5408503") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F9"), "https://www.lancers.jp/work/detail/5400689") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F10"), "https://www.lancers.jp/work/detail/5393695") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F11"), "https://www.lancers.jp/work/detail/5408563") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F12"), "https://www.lancers.jp/work/detail/5408286") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F13"), "https://www.lancers.jp/work/detail/5398382") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F14"), "https://www.lancers.jp/work/detail/5408148") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F15"), "https://www.lancers.jp/work/detail/5408554") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F16"), "https://www.lancers.jp/work/detail/5408546") | Out-Null
